# Auto commit at 2025-10-22  7:59:08.96
# Refresh the cumulative metrics on the "Metrics" sheet with the latest
# counters. Everything downstream (the "today" sheet's Metrics! lookups,
# running totals in column E/F, etc.) recalculates automatically because
# those cells hold formulas, not hard-coded numbers.

$wb = $excel.ActiveWorkbook

$metrics = $wb.Worksheets.Item("Metrics")
$metrics.Range("B2").Value  = 290791.73000000004
$metrics.Range("B3").Value  = 238227.04
$metrics.Range("B4").Value  = 92648.89999999998
$metrics.Range("B5").Value  = 11568
$metrics.Range("B6").Value  = 4657923.2
$metrics.Range("B7").Value  = 3928045.7099999995
$metrics.Range("B8").Value  = 1363251.04
$metrics.Range("B9").Value  = 180569
$metrics.Range("B10").Value = 33123247.000999823
$metrics.Range("B11").Value = 31203267.23
$metrics.Range("B12").Value = 11644959.930000002
$metrics.Range("B13").Value = 1278196

# Leave the cursor where the author left it when the workbook was saved.
$metrics.Range("E22").Select()

$today = $wb.Worksheets.Item("today")
$today.Activate()
$today.Range("F7").Select()
